$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P6").Value = 274250
$ws.Range("P7").Value = 110388
$ws.Range("P8").Value = 303090
$ws.Range("P9").Value = 310986
$ws.Range("P10").Value = 311727
$ws.Range("P11").Value = 284604
$ws.Range("P12").Value = 441228
$ws.Range("P13").Value = 282884
